# Update 想去人数 (number of people interested) figures that changed
# between data refreshes for the gh-pages generated output.
#
# Sheet "展览" (Exhibitions): rows 2-5, column F
# Sheet "全部类型" (All types): rows 2,3,6,7, column F

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 102
$wsExhibition.Range("F3").Value = 2180
$wsExhibition.Range("F4").Value = 894
$wsExhibition.Range("F5").Value = 1536

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 102
$wsAll.Range("F3").Value = 2180
$wsAll.Range("F6").Value = 894
$wsAll.Range("F7").Value = 1536
